$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '320.95'
Set-TextValue $ws.Range('E2') '6.36%'
Set-TextValue $ws.Range('D3') '49.32'
Set-TextValue $ws.Range('E3') '12.26%'
Set-TextValue $ws.Range('D4') '5.340'
Set-TextValue $ws.Range('E4') '5.01%'
Set-TextValue $ws.Range('D5') '0.08060'
Set-TextValue $ws.Range('E5') '4.72%'
Set-TextValue $ws.Range('D6') '4.611'
Set-TextValue $ws.Range('E6') '4.26%'
Set-TextValue $ws.Range('D7') '1.381'
Set-TextValue $ws.Range('E7') '32.65%'
Set-TextValue $ws.Range('D8') '1.640'
Set-TextValue $ws.Range('E8') '1.54%'
Set-TextValue $ws.Range('D9') '0.1257'
Set-TextValue $ws.Range('E9') '-2.23%'
Set-TextValue $ws.Range('E10') '5.45%'
Set-TextValue $ws.Range('D11') '0.09729'
Set-TextValue $ws.Range('E11') '5.05%'
Set-TextValue $ws.Range('D12') '0.04717'
Set-TextValue $ws.Range('E12') '12.35%'
Set-TextValue $ws.Range('D13') '0.1046'
Set-TextValue $ws.Range('E13') '-0.13%'
Set-TextValue $ws.Range('D14') '0.001315'
Set-TextValue $ws.Range('E14') '2.76%'
Set-TextValue $ws.Range('D15') '0.04200'
Set-TextValue $ws.Range('E15') '0.44%'
Set-TextValue $ws.Range('D16') '0.005799'
Set-TextValue $ws.Range('E16') '0.63%'
Set-TextValue $ws.Range('B17') 'HotbitToken'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range('D17') '0.004326'
Set-TextValue $ws.Range('E17') '-2.02%'
Set-TextValue $ws.Range('B18') 'LEO'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D18') '3.347'
Set-TextValue $ws.Range('E18') '0.06%'
Set-TextValue $ws.Range('B19') 'BTSEToken'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D19') '2.455'
Set-TextValue $ws.Range('E19') '5.36%'
Set-TextValue $ws.Range('B20') 'BitpandaEcosystemToken'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range('D20') '0.3506'
Set-TextValue $ws.Range('E20') '4.66%'
Set-TextValue $ws.Range('B21') 'MCDex'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D21') '8.024'
Set-TextValue $ws.Range('E21') '-4.01%'
Set-TextValue $ws.Range('B22') 'ProBitToken'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range('D22') '0.1361'
Set-TextValue $ws.Range('E22') '-2.74%'
Set-TextValue $ws.Range('B23') 'ZBToken'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range('D23') '0.3093'
Set-TextValue $ws.Range('E23') '-2.63%'
Set-TextValue $ws.Range('B24') 'BitKan'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range('D24') '0.001315'
Set-TextValue $ws.Range('E24') '2.38%'
Set-TextValue $ws.Range('D25') '0.0001349'
Set-TextValue $ws.Range('E25') '-0.04%'
Set-TextValue $ws.Range('D26') '0.0003537'
Set-TextValue $ws.Range('E26') '-95.28%'
Set-TextValue $ws.Range('D38') '0.02739'
Set-TextValue $ws.Range('E38') '10.14%'
Set-TextValue $ws.Range('D39') '0.06098'
Set-TextValue $ws.Range('E39') '15.14%'
Set-TextValue $ws.Range('D40') '0.01086'
Set-TextValue $ws.Range('E40') '83.01%'
Set-TextValue $ws.Range('D41') '0.008048'
Set-TextValue $ws.Range('E41') '4.47%'
Set-TextValue $ws.Range('D42') '0.1464'
Set-TextValue $ws.Range('E42') '8.67%'
Set-TextValue $ws.Range('D43') '0.007900'
Set-TextValue $ws.Range('E43') '7.25%'
Set-TextValue $ws.Range('D44') '0.008666'
Set-TextValue $ws.Range('E44') '14.73%'
Set-TextValue $ws.Range('D45') '0.3508'
Set-TextValue $ws.Range('E45') '16.91%'
Set-TextValue $ws.Range('D46') '0.00006849'
Set-TextValue $ws.Range('E46') '2.66%'
Set-TextValue $ws.Range('D47') '0.00000000750'
Set-TextValue $ws.Range('E47') '0.01%'
Set-TextValue $ws.Range('D48') '0.05922'
Set-TextValue $ws.Range('E48') '39.65%'
Set-TextValue $ws.Range('D49') '0.003998'
Set-TextValue $ws.Range('E49') '-4.83%'
Set-TextValue $ws.Range('D50') '0.00002099'
Set-TextValue $ws.Range('E50') '0.01%'
Set-TextValue $ws.Range('D51') '0.0001999'
Set-TextValue $ws.Range('E51') '0.01%'
